# "Weapons" is the first / active sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weapons")

# Insert a new row 14 ("Dremora" weapon material), pushing the existing
# rows 14-38 down to 15-39.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "Dremora"
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 10

# Match the author's final cursor/selection position and scroll state.
$ws.Activate()
$ws.Range("H14").Select()
